$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C6").Value2 = 63801.39240103823
$ws.Range("C7").Value2 = 57421.25316093441
$ws.Range("C9").Value2 = 11804.219129958965
$ws.Range("C10").Value2 = 55637.17327107927
$ws.Range("C11").Value2 = 51997.17327107927
$ws.Range("C12").Value2 = 17640.0
$ws.Range("C13").Value2 = 14000.0
$ws.Range("C14").Value2 = 37997.17327107927
$ws.Range("C15").Value2 = 37219.05728290892
$ws.Range("C16").Value2 = 319.0286971703372
$ws.Range("C18").Value2 = 2095.8
$ws.Range("C19").Value2 = 35442.28598007928
$ws.Range("C20").Value2 = 20733.38979152061

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C2").Value2 = 6380.573943406745
$ws.Range("C3").Value2 = 8157.333333333332
$ws.Range("D3").Value2 = 27.846388204035645
$ws.Range("C5").Value2 = 8157.333333333332
$ws.Range("C8").Value2 = 10243.0
$ws.Range("D8").Value2 = 60.534147724820706
$ws.Range("D9").Value2 = 0.24176597168341904
$ws.Range("C10").Value2 = 7149.0
$ws.Range("D10").Value2 = 12.043212153152714
$ws.Range("C11").Value2 = 14388.0
$ws.Range("D11").Value2 = 125.4969557224173
$ws.Range("C12").Value2 = 6463.0
$ws.Range("D12").Value2 = 1.291828248122254
$ws.Range("C13").Value2 = 7891.0
$ws.Range("D13").Value2 = 23.67226005043056
$ws.Range("C14").Value2 = 10802.0
$ws.Range("D14").Value2 = 69.29511507600442
$ws.Range("C15").Value2 = 21031.0
$ws.Range("D15").Value2 = 229.6098468027633

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value2 = 6763.408380011148
$ws.Range("C3").Value2 = 7181.0
$ws.Range("D3").Value2 = 6.17427776833672
$ws.Range("C5").Value2 = 7180.999999999999
$ws.Range("A8").Value2 = "KROO"
$ws.Range("C8").Value2 = 7561.0
$ws.Range("D8").Value2 = 11.792746721402883
$ws.Range("A9").Value2 = "TORENBEEK_1982"
$ws.Range("C9").Value2 = 6631.0
$ws.Range("D9").Value2 = -1.9577167689958312
$ws.Range("A10").Value2 = "RAYMER"
$ws.Range("C10").Value2 = 8394.0
$ws.Range("D10").Value2 = 24.1090220843084
$ws.Range("A11").Value2 = "TORENBEEK_2013"
$ws.Range("C11").Value2 = 6138.0
$ws.Range("D11").Value2 = -9.246940963368484

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C2").Value2 = 733.7660034917756
$ws.Range("C3").Value2 = 796.625
$ws.Range("D3").Value2 = 8.56662699131563
$ws.Range("C5").Value2 = 796.6249999999999
$ws.Range("A8").Value2 = "SADRAEY"
$ws.Range("C8").Value2 = 1040.0
$ws.Range("D8").Value2 = 41.734557754236
$ws.Range("D9").Value2 = 92.84076848292686
$ws.Range("A10").Value2 = "KROO"
$ws.Range("C10").Value2 = 737.0
$ws.Range("D10").Value2 = 0.4407394854537826
$ws.Range("C11").Value2 = 399.0
$ws.Range("D11").Value2 = -45.62299178467292
$ws.Range("A12").Value2 = "ROSKAM"
$ws.Range("C12").Value2 = 1523.0
$ws.Range("D12").Value2 = 107.55935717278983
$ws.Range("A13").Value2 = "RAYMER"
$ws.Range("C13").Value2 = 507.0
$ws.Range("D13").Value2 = -30.904403094809947
$ws.Range("A14").Value2 = "TORENBEEK_1976"
$ws.Range("C14").Value2 = 52.0
$ws.Range("D14").Value2 = -92.91327211228821
$ws.Range("C15").Value2 = 700.0
$ws.Range("D15").Value2 = -4.601739973110383

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value2 = 733.7660034917756
$ws.Range("C3").Value2 = 673.0
$ws.Range("D3").Value2 = -8.281387145576124
$ws.Range("C5").Value2 = 672.9999999999999
$ws.Range("A8").Value2 = "SADRAEY"
$ws.Range("C8").Value2 = 749.0
$ws.Range("D8").Value2 = 2.07613822877189
$ws.Range("D9").Value2 = 56.04429675826945
$ws.Range("A10").Value2 = "KROO"
$ws.Range("C10").Value2 = 488.0
$ws.Range("D10").Value2 = -33.49378443839695
$ws.Range("A11").Value2 = "ROSKAM"
$ws.Range("C11").Value2 = 1523.0
$ws.Range("D11").Value2 = 107.55935717278983
$ws.Range("A12").Value2 = "RAYMER"
$ws.Range("C12").Value2 = 180.0
$ws.Range("D12").Value2 = -75.46901885022838
$ws.Range("A13").Value2 = "TORENBEEK_1976"
$ws.Range("C13").Value2 = 124.0
$ws.Range("D13").Value2 = -83.10087965237955
$ws.Range("C14").Value2 = 502.0
$ws.Range("D14").Value2 = -31.58581923785916

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C2").Value2 = 1212.309049247281
$ws.Range("C3").Value2 = 1390.0
$ws.Range("D3").Value2 = 14.657232069912094
$ws.Range("C5").Value2 = 1389.9999999999998
$ws.Range("A10").Value2 = "ROSKAM"
$ws.Range("C10").Value2 = 687.0
$ws.Range("D10").Value2 = 13.337436592848357
$ws.Range("D11").Value2 = 14.492257635279127
$ws.Range("C12").Value2 = 704.0
$ws.Range("D12").Value2 = 16.142001981608797
$ws.Range("C14").Value2 = 695.0
$ws.Range("A17").Value2 = "ROSKAM"
$ws.Range("C17").Value2 = 687.0
$ws.Range("D17").Value2 = 13.337436592848357
$ws.Range("D18").Value2 = 14.492257635279127
$ws.Range("C19").Value2 = 704.0
$ws.Range("D19").Value2 = 16.142001981608797
$ws.Range("C21").Value2 = 695.0

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C2").Value2 = 5295.876373027598
$ws.Range("C3").Value2 = 6450.666666666666
$ws.Range("D3").Value2 = 21.80546168941039
$ws.Range("C5").Value2 = 6450.666666666665
$ws.Range("A11").Value2 = "TORENBEEK_2013"
$ws.Range("C11").Value2 = 3457.0
$ws.Range("D11").Value2 = 30.55440710839968
$ws.Range("A12").Value2 = "TORENBEEK_1976"
$ws.Range("C12").Value2 = 2954.0
$ws.Range("D12").Value2 = 11.558495400119368
$ws.Range("A13").Value2 = "KUNDU"
$ws.Range("C13").Value2 = 3265.0
$ws.Range("D13").Value2 = 23.303482559712165
$ws.Range("C14").Value2 = 3225.333333333333
$ws.Range("A18").Value2 = "TORENBEEK_2013"
$ws.Range("C18").Value2 = 3457.0
$ws.Range("D18").Value2 = 30.55440710839968
$ws.Range("A19").Value2 = "TORENBEEK_1976"
$ws.Range("C19").Value2 = 2954.0
$ws.Range("D19").Value2 = 11.558495400119368
$ws.Range("A20").Value2 = "KUNDU"
$ws.Range("C20").Value2 = 3265.0
$ws.Range("D20").Value2 = 23.303482559712165
$ws.Range("C21").Value2 = 3225.333333333333

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value2 = 2616.035316796765
$ws.Range("C3").Value2 = 2535.4314581872886
$ws.Range("D3").Value2 = -3.0811456593090853
$ws.Range("C5").Value2 = 2535.4314581872886
$ws.Range("A9").Value2 = "TORENBEEK_1976"
$ws.Range("C9").Value2 = 2535.431458187289
$ws.Range("D9").Value2 = -3.0811456593090774
$ws.Range("A11").Value2 = "TORENBEEK_1976"
$ws.Range("C11").Value2 = 391.4126375939395
$ws.Range("A13").Value2 = "TORENBEEK_1976"
$ws.Range("C13").Value2 = 2144.01882059335

$ws = $wb.Worksheets.Item("SYSTEMS")
$ws.Range("C2").Value2 = 8677.580563033172
$ws.Range("C3").Value2 = 8258.229521892012
$ws.Range("D3").Value2 = -4.832580211673437
$ws.Range("C4").Value2 = 8258.22952189201
$ws.Range("C8").Value2 = 8258.229521892012
$ws.Range("D8").Value2 = -4.832580211673409
$ws.Range("C11").Value2 = 336.84274966573867
$ws.Range("C13").Value2 = 336.8427496657386
$ws.Range("C21").Value2 = 1033.7492545113123
$ws.Range("C23").Value2 = 1033.749254511312
$ws.Range("C26").Value2 = 531.6667122518941
$ws.Range("C28").Value2 = 531.666712251894
$ws.Range("C36").Value2 = 785.0973560969437
$ws.Range("C38").Value2 = 785.0973560969436
$ws.Range("C41").Value2 = 3321.8094918179895
$ws.Range("C43").Value2 = 3321.8094918179895
